$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) updates -- stored as text in the sheet, so prefix with
# an apostrophe to force a text literal even though the content looks numeric.
$ws.Range("D2").Value  = "'268.48"
$ws.Range("D3").Value  = "'22.87"
$ws.Range("D4").Value  = "'6.332"
$ws.Range("D5").Value  = "'0.06174"
$ws.Range("D7").Value  = "'6.664"
$ws.Range("D9").Value  = "'0.8294"
$ws.Range("D10").Value = "'0.01373"
$ws.Range("D11").Value = "'0.1607"
$ws.Range("D12").Value = "'0.08307"
$ws.Range("D13").Value = "'0.03475"
$ws.Range("D14").Value = "'0.03193"
$ws.Range("D15").Value = "'0.09328"
$ws.Range("D16").Value = "'3.841"
$ws.Range("D17").Value = "'0.001653"
$ws.Range("D18").Value = "'0.04746"
$ws.Range("D19").Value = "'0.006329"
$ws.Range("D20").Value = "'0.005659"
$ws.Range("D23").Value = "'3.727"
$ws.Range("D25").Value = "'0.3304"
$ws.Range("D26").Value = "'0.1238"
$ws.Range("D40").Value = "'0.04715"
$ws.Range("D41").Value = "'0.006976"

# Rows 42/43 swap places (BKEXToken <-> CEJI) along with their data.
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.003801"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "'0.1158"
$ws.Range("E43").Value = "42BKEXTokenBKK"

$ws.Range("D44").Value = "'0.01178"
$ws.Range("D45").Value = "'0.00006250"
$ws.Range("D46").Value = "'0.0009902"
$ws.Range("D48").Value = "'0.9202"
$ws.Range("D49").Value = "'0.002329"
$ws.Range("D50").Value = "'0.00001400"
$ws.Range("E50").Value = "49CryptobidCoinCBCWorstin24h"
